$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RunToTest" (A) / "ConsignmentID" (B) block for rows 22-29 is being
# re-shuffled: the "Y" (run) flag moves from rows 22-23 down to rows 24-25,
# and the ConsignmentID values are bumped forward by one slot, with a brand
# new UAT42092150 value introduced at the bottom (row 29).

$ws.Range("A22").Value = "N"
$ws.Range("B22").Value = "UAT42092143"

$ws.Range("A23").Value = "N"
$ws.Range("B23").Value = "UAT42092144"

$ws.Range("A24").Value = "Y"
$ws.Range("B24").Value = "UAT42092145"

$ws.Range("A25").Value = "Y"
$ws.Range("B25").Value = "UAT42092146"

$ws.Range("A26").Value = "N"
$ws.Range("B26").Value = "UAT42092147"

$ws.Range("A27").Value = "N"
$ws.Range("B27").Value = "UAT42092148"

$ws.Range("A28").Value = "N"
$ws.Range("B28").Value = "UAT42092149"

$ws.Range("A29").Value = "N"
$ws.Range("B29").Value = "UAT42092150"

# Restore the window's view/selection state to match the saved workbook.
$ws.Range("B34").Select() | Out-Null
